$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.596
$ws.Range("D2").Value = 0.172
$ws.Range("E2").Value = -0.016
$ws.Range("F2").Value = -0.119
$ws.Range("G2").Value = 0.067
$ws.Range("H2").Value = -0.111
$ws.Range("I2").Value = 0.152
$ws.Range("J2").Value = 0.14
$ws.Range("K2").Value = -0.127
$ws.Range("L2").Value = 0.014
$ws.Range("B3").Value = 0.596
$ws.Range("D3").Value = 0.197
$ws.Range("E3").Value = -0.276
$ws.Range("F3").Value = -0.352
$ws.Range("G3").Value = -0.012
$ws.Range("H3").Value = 0.047
$ws.Range("I3").Value = -0.018
$ws.Range("J3").Value = -0.047
$ws.Range("K3").Value = -0.092
$ws.Range("L3").Value = 0.137
$ws.Range("B4").Value = 0.172
$ws.Range("C4").Value = 0.197
$ws.Range("E4").Value = 0.479
$ws.Range("F4").Value = 0.073
$ws.Range("G4").Value = 0.622
$ws.Range("H4").Value = 0.318
$ws.Range("I4").Value = 0.423
$ws.Range("J4").Value = 0.357
$ws.Range("K4").Value = 0.095
$ws.Range("L4").Value = -0.067
$ws.Range("B5").Value = -0.016
$ws.Range("C5").Value = -0.276
$ws.Range("D5").Value = 0.479
$ws.Range("F5").Value = 0.646
$ws.Range("G5").Value = 0.626
$ws.Range("H5").Value = 0.339
$ws.Range("I5").Value = 0.484
$ws.Range("J5").Value = 0.423
$ws.Range("K5").Value = 0.158
$ws.Range("L5").Value = -0.162
$ws.Range("B6").Value = -0.119
$ws.Range("C6").Value = -0.352
$ws.Range("D6").Value = 0.073
$ws.Range("E6").Value = 0.646
$ws.Range("G6").Value = 0.257
$ws.Range("H6").Value = 0.218
$ws.Range("I6").Value = 0.297
$ws.Range("J6").Value = 0.209
$ws.Range("K6").Value = 0.213
$ws.Range("L6").Value = -0.287
$ws.Range("B7").Value = 0.067
$ws.Range("C7").Value = -0.012
$ws.Range("D7").Value = 0.622
$ws.Range("E7").Value = 0.626
$ws.Range("F7").Value = 0.257
$ws.Range("H7").Value = 0.381
$ws.Range("I7").Value = 0.354
$ws.Range("J7").Value = 0.312
$ws.Range("K7").Value = 0.016
$ws.Range("L7").Value = 0.029
$ws.Range("B8").Value = -0.111
$ws.Range("C8").Value = 0.047
$ws.Range("D8").Value = 0.318
$ws.Range("E8").Value = 0.339
$ws.Range("F8").Value = 0.218
$ws.Range("G8").Value = 0.381
$ws.Range("I8").Value = 0.037
$ws.Range("J8").Value = 0.115
$ws.Range("K8").Value = 0.091
$ws.Range("L8").Value = -0.02
$ws.Range("B9").Value = 0.152
$ws.Range("C9").Value = -0.018
$ws.Range("D9").Value = 0.423
$ws.Range("E9").Value = 0.484
$ws.Range("F9").Value = 0.297
$ws.Range("G9").Value = 0.354
$ws.Range("H9").Value = 0.037
$ws.Range("J9").Value = 0.601
$ws.Range("L9").Value = -0.31
$ws.Range("B10").Value = 0.14
$ws.Range("C10").Value = -0.047
$ws.Range("D10").Value = 0.357
$ws.Range("E10").Value = 0.423
$ws.Range("F10").Value = 0.209
$ws.Range("G10").Value = 0.312
$ws.Range("H10").Value = 0.115
$ws.Range("I10").Value = 0.601
$ws.Range("K10").Value = -0.099
$ws.Range("L10").Value = -0.222
$ws.Range("B11").Value = -0.127
$ws.Range("C11").Value = -0.092
$ws.Range("D11").Value = 0.095
$ws.Range("E11").Value = 0.158
$ws.Range("F11").Value = 0.213
$ws.Range("G11").Value = 0.016
$ws.Range("H11").Value = 0.091
$ws.Range("J11").Value = -0.099
$ws.Range("L11").Value = -0.047
$ws.Range("B12").Value = 0.014
$ws.Range("C12").Value = 0.137
$ws.Range("D12").Value = -0.067
$ws.Range("E12").Value = -0.162
$ws.Range("F12").Value = -0.287
$ws.Range("G12").Value = 0.029
$ws.Range("H12").Value = -0.02
$ws.Range("I12").Value = -0.31
$ws.Range("J12").Value = -0.222
$ws.Range("K12").Value = -0.047
